# Weekly data refresh: insert two fresh rows of data (week of 44714) at the
# top of the "Poroto verde" block (rows 486-487), pushing the existing rows
# 486-522 down to 488-524.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 486:522 down by two rows, creating blank rows 486:487.
$ws.Range("A486:R487").Insert()

# New row 486
$ws.Cells.Item(486, 1).Value = 9
$ws.Cells.Item(486, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(486, 3).Value = "Metropolitana"
$ws.Cells.Item(486, 4).Value = 44714
$ws.Cells.Item(486, 5).Value = 13
$ws.Cells.Item(486, 6).Value = 100112031
$ws.Cells.Item(486, 7).Value = "Poroto verde"
$ws.Cells.Item(486, 8).Value = "Magnum"
$ws.Cells.Item(486, 9).Value = "Primera"
$ws.Cells.Item(486, 10).Value = 61
$ws.Cells.Item(486, 11).Value = 22000
$ws.Cells.Item(486, 12).Value = 23000
$ws.Cells.Item(486, 13).Value = 22508
$ws.Cells.Item(486, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(486, 15).Value = "Perú"
$ws.Cells.Item(486, 16).Value = 900
$ws.Cells.Item(486, 17).Value = 25
$ws.Cells.Item(486, 18).Value = "Hortaliza"

# New row 487
$ws.Cells.Item(487, 1).Value = 9
$ws.Cells.Item(487, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(487, 3).Value = "Metropolitana"
$ws.Cells.Item(487, 4).Value = 44714
$ws.Cells.Item(487, 5).Value = 13
$ws.Cells.Item(487, 6).Value = 100112031
$ws.Cells.Item(487, 7).Value = "Poroto verde"
$ws.Cells.Item(487, 8).Value = "Magnum"
$ws.Cells.Item(487, 9).Value = "Primera"
$ws.Cells.Item(487, 10).Value = 43
$ws.Cells.Item(487, 11).Value = 23000
$ws.Cells.Item(487, 12).Value = 24000
$ws.Cells.Item(487, 13).Value = 23488
$ws.Cells.Item(487, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(487, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(487, 16).Value = 940
$ws.Cells.Item(487, 17).Value = 25
$ws.Cells.Item(487, 18).Value = "Hortaliza"
